# Auto-generated edit script: updates profit-tracking cells across all 8 sheets
# per the authoritative diff (scheduled market-price refresh).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H68").Value = 89999
$ws.Range("J68").Value = 89999
$ws.Range("L68").Value = 89999
$ws.Range("N68").Value = -91497
$ws.Range("H71").Value = 89999
$ws.Range("J71").Value = 89999
$ws.Range("L71").Value = 269997
$ws.Range("N71").Value = -277485
$ws.Range("H74").Value = 6915.8
$ws.Range("I74").Value = 6915.8
$ws.Range("K74").Value = 6915.8
$ws.Range("M74").Value = -5979.8
$ws.Range("H77").Value = 6915.8
$ws.Range("I77").Value = 6915.8
$ws.Range("K77").Value = 34579
$ws.Range("M77").Value = -29899
$ws.Range("H95").Value = 34995
$ws.Range("J95").Value = 34995
$ws.Range("L95").Value = 34995
$ws.Range("N95").Value = -40487

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4304
$ws.Range("I61").Value = 4304
$ws.Range("K61").Value = 4304
$ws.Range("M61").Value = -4092
$ws.Range("H136").Value = 4304
$ws.Range("I136").Value = 4304
$ws.Range("K136").Value = 12912
$ws.Range("M136").Value = -10362
$ws.Range("H139").Value = 73678.75
$ws.Range("J139").Value = 73678.75
$ws.Range("L139").Value = 73678.75
$ws.Range("N139").Value = -83958.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 13000
$ws.Range("J76").Value = 13000
$ws.Range("L76").Value = 13000
$ws.Range("N76").Value = -13630
$ws.Range("H79").Value = 13000
$ws.Range("J79").Value = 13000
$ws.Range("L79").Value = 13000
$ws.Range("N79").Value = -15184
$ws.Range("H88").Value = 45990
$ws.Range("J88").Value = 45990
$ws.Range("L88").Value = 45990
$ws.Range("N88").Value = -46802
$ws.Range("H91").Value = 45990
$ws.Range("J91").Value = 45990
$ws.Range("L91").Value = 45990
$ws.Range("N91").Value = -48798
$ws.Range("H95").Value = 19980
$ws.Range("J95").Value = 19980
$ws.Range("L95").Value = 19980
$ws.Range("N95").Value = -25472

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 18152
$ws.Range("J43").Value = 18152
$ws.Range("L43").Value = 18152
$ws.Range("N43").Value = -18520
$ws.Range("H74").Value = 56704.715
$ws.Range("J74").Value = 56704.715
$ws.Range("L74").Value = 56704.715
$ws.Range("N74").Value = -58452.715
$ws.Range("H77").Value = 56704.715
$ws.Range("J77").Value = 56704.715
$ws.Range("L77").Value = 170114.145
$ws.Range("N77").Value = -178850.145
$ws.Range("H101").Value = 18152
$ws.Range("J101").Value = 18152
$ws.Range("L101").Value = 18152
$ws.Range("N101").Value = -24642
$ws.Range("H132").Value = 3742.2
$ws.Range("I132").Value = 2801.4285
$ws.Range("J132").Value = 5937.3335
$ws.Range("K132").Value = 8404.2855
$ws.Range("L132").Value = 17812.0005
$ws.Range("M132").Value = -5874.2855
$ws.Range("N132").Value = -22872.0005
$ws.Range("H134").Value = 5723.222
$ws.Range("I134").Value = 5723.222
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 17169.666
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -14634.666
$ws.Range("N134").Value = $null

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 430
$ws.Range("I14").Value = 430
$ws.Range("K14").Value = 1290
$ws.Range("M14").Value = -1117
$ws.Range("H69").Value = 7000
$ws.Range("I69").Value = 1000
$ws.Range("J69").Value = 10000
$ws.Range("K69").Value = 3000
$ws.Range("L69").Value = 30000
$ws.Range("M69").Value = -2189
$ws.Range("N69").Value = -31622
$ws.Range("H72").Value = 7000
$ws.Range("I72").Value = 1000
$ws.Range("J72").Value = 10000
$ws.Range("K72").Value = 9000
$ws.Range("L72").Value = 90000
$ws.Range("M72").Value = -4944
$ws.Range("N72").Value = -98112
$ws.Range("H129").Value = 1266.6666
$ws.Range("I129").Value = 1266.6666
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 3799.9998
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 1200.0002
$ws.Range("N129").Value = $null

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 20000
$ws.Range("I43").Value = 20000
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 20000
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -19849
$ws.Range("N43").Value = $null
$ws.Range("H102").Value = 1849.8334
$ws.Range("I102").Value = 1739.8
$ws.Range("K102").Value = 1739.8
$ws.Range("M102").Value = -117.8
$ws.Range("H132").Value = 3708
$ws.Range("I132").Value = 3315.6667
$ws.Range("K132").Value = 9947.000100000001
$ws.Range("M132").Value = -7417.000100000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2125.375
$ws.Range("I7").Value = 2071.8572
$ws.Range("K7").Value = 2071.8572
$ws.Range("M7").Value = -1959.8572
$ws.Range("H10").Value = 4999
$ws.Range("I10").Value = 4999
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 4999
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -4859
$ws.Range("N10").Value = $null
$ws.Range("H22").Value = 11692.571
$ws.Range("J22").Value = 8714.286
$ws.Range("L22").Value = 8714.286
$ws.Range("N22").Value = -9304.286
$ws.Range("H27").Value = 11692.571
$ws.Range("J27").Value = 8714.286
$ws.Range("L27").Value = 8714.286
$ws.Range("N27").Value = -8928.286
$ws.Range("H126").Value = 2125.375
$ws.Range("I126").Value = 2071.8572
$ws.Range("K126").Value = 6215.571599999999
$ws.Range("M126").Value = -3745.571599999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 89999
$ws.Range("J70").Value = 89999
$ws.Range("L70").Value = 89999
$ws.Range("N70").Value = -90629
$ws.Range("H73").Value = 89999
$ws.Range("J73").Value = 89999
$ws.Range("L73").Value = 89999
$ws.Range("N73").Value = -92183
$ws.Range("H81").Value = 3628.1333
$ws.Range("I81").Value = 951.8333
$ws.Range("K81").Value = 1903.6666
$ws.Range("M81").Value = -842.6666
$ws.Range("H84").Value = 3628.1333
$ws.Range("I84").Value = 951.8333
$ws.Range("K84").Value = 9518.333000000001
$ws.Range("M84").Value = -4214.333000000001
$ws.Range("H98").Value = 29999
$ws.Range("J98").Value = 29999
$ws.Range("L98").Value = 29999
$ws.Range("N98").Value = -35989
$ws.Range("H113").Value = 1241.3334
$ws.Range("I113").Value = 1371
$ws.Range("K113").Value = 4113
$ws.Range("M113").Value = -1943
$ws.Range("H136").Value = 7980
$ws.Range("I136").Value = 8725
$ws.Range("K136").Value = 26175
$ws.Range("M136").Value = -23625

